$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.047.30"
$ws.Range("E2").Value = "  -1.14%  "
$ws.Range("D3").Value = "1.642.95"
$ws.Range("E3").Value = "  -1.50%  "
$ws.Range("E4").Value = "  -0.64%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.44"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.92%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5182"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.10%  "
$ws.Range("E7").Value = "  -0.57%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2615"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06275"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.06%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.41"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.24%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07739"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.35%  "
$ws.Range("D12").Value = "1.708.99"
$ws.Range("E12").Value = "  +2.52%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.470"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.22%  "
$ws.Range("D14").Value = "1.868.36"
$ws.Range("E14").Value = "  -1.46%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5568"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.48%  "
$ws.Range("D16").Value = "0.0₅7985"
$ws.Range("E16").Value = "  -2.56%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.67"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.79%  "
$ws.Range("D18").Value = "26.042.72"
$ws.Range("E18").Value = "  -1.25%  "
$ws.Range("E19").Value = "  -0.56%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.618"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.70%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "192.58"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.59%  "
$ws.Range("E22").Value = "  -2.43%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.947"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.79%  "
$ws.Range("E24").Value = "  -0.61%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.63"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.21%  "
$ws.Range("E26").Value = "  -2.73%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.152"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.17%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.88"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.52%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.477"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.60%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05641"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.86%  "
$ws.Range("E31").Value = "  -1.85%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.452"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.19%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.355"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.99%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.594"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.95%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.788"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.36%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.410"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.45%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9373"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.56%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5652"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.32%  "
$ws.Range("E39").Value = "  +1.77%  "
$ws.Range("E40").Value = "  -1.94%  "
$ws.Range("D41").Value = "1.051.06"
$ws.Range("E41").Value = "  -1.27%  "
$ws.Range("E42").Value = "  -0.63%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8428"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.06%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "102.31"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.60%  "
$ws.Range("D45").Value = "1.779.71"
$ws.Range("E45").Value = "  -1.50%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "56.83"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.91%  "
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "0.0₈105"
$ws.Range("E47").Value = "  -0.96%  "
$ws.Range("B48").Value = "Frax"
$ws.Range("C48").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.008"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.69%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05307"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.68%  "
$ws.Range("E50").Value = "  -1.43%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.910"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.33%  "
